# Update Investigation for correctness
# - Replace the investigation identifier value (was "testARC_licenseMissing")
#   with "testARC_correctOrcid" across B7:B9 (date-ish style carried over from
#   the source workbook).
# - Add three contacts (Person X / Y / Z) with emails, affiliation links and
#   a new "Comment[ORCID]" row with ORCID links for Person X.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Investigation Identifier / Title / Description get the corrected value ---
$ws.Range("B7").Value = "testARC_correctOrcid"
$ws.Range("B7").NumberFormat = "d-mmm-yy"

$ws.Range("B8").Value = "testARC_correctOrcid"
$ws.Range("B8").NumberFormat = "d-mmm-yy"

$ws.Range("B9").Value = "testARC_correctOrcid"
$ws.Range("B9").NumberFormat = "d-mmm-yy"

# --- Investigation Person Last / First Name / Affiliation for 3 contacts ---
$ws.Range("B21").Value = "Person X"
$ws.Range("C21").Value = "Person Y"
$ws.Range("D21").Value = "Person Z"

$ws.Range("B22").Value = "Person X"
$ws.Range("C22").Value = "Person Y"
$ws.Range("D22").Value = "Person Z"

$ws.Range("B28").Value = "Person X"
$ws.Range("C28").Value = "Person Y"
$ws.Range("D28").Value = "Person Z"

# --- Investigation Person Email, with mailto hyperlinks ---
$ws.Hyperlinks.Add($ws.Range("B24"), "mailto:per@son.X", "", "", "per@son.X") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C24"), "mailto:per@son.Y", "", "", "per@son.Y") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D24"), "mailto:per@son.Z", "", "", "per@son.Z") | Out-Null

# --- New row 32: Comment[ORCID] with http/https ORCID links for Person X ---
$ws.Range("A32").Value = "Comment[ORCID]"
$ws.Hyperlinks.Add($ws.Range("B32"), "http://orcid.org/0000-0001-5109-3700", "", "", "http://orcid.org/0000-0001-5109-3700") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C32"), "https://orcid.org/0000-0001-5109-3700", "", "", "https://orcid.org/0000-0001-5109-3700") | Out-Null
$ws.Range("D32").Value = "0000-0001-5109-3700"

# --- Column widths to fit the new contact / link columns ---
# (values chosen so the saved/quantized OOXML column width lands as close as
# possible to the target widths of 51.71 / 34.43 / 35.29 / 19.29 chars)
$ws.Columns.Item(1).ColumnWidth = 50.75
$ws.Columns.Item(2).ColumnWidth = 33.59
$ws.Columns.Item(3).ColumnWidth = 34.42
$ws.Columns.Item(4).ColumnWidth = 18.42

# --- Keep active selection sane, matching the saved view of the workbook ---
$ws.Range("C32").Select() | Out-Null
